$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllEntries")
$ws.Columns("AN:AO").Insert()
